$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 293 (pushes old rows 293:335 down to 295:337,
# inheriting formatting/styles from the rows immediately above, matching the
# existing data rows' look - date column D keeps its date style).
$ws.Rows("293:294").Insert()

# Fill in the new row 293 with its data.
$ws.Range("A293").Value = 4
$ws.Range("B293").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C293").Value = "Los Lagos"
$ws.Range("D293").Value = 44505
$ws.Range("E293").Value = 10
$ws.Range("F293").Value = 100112004
$ws.Range("G293").Value = "Cebolla"
$ws.Range("H293").Value = "Morada(o)"
$ws.Range("I293").Value = "1a nueva(o)"
$ws.Range("J293").Value = 140
$ws.Range("K293").Value = 10000
$ws.Range("L293").Value = 11000
$ws.Range("M293").Value = 10500
$ws.Range("N293").Value = "$/malla 18 kilos"
$ws.Range("O293").Value = "Región de Arica y Parinacota"
$ws.Range("P293").Value = 583
$ws.Range("Q293").Value = 18
$ws.Range("R293").Value = "Hortaliza"

# Fill in the new row 294 with its data.
$ws.Range("A294").Value = 4
$ws.Range("B294").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C294").Value = "Los Lagos"
$ws.Range("D294").Value = 44505
$ws.Range("E294").Value = 10
$ws.Range("F294").Value = 100112004
$ws.Range("G294").Value = "Cebolla"
$ws.Range("H294").Value = "Sin especificar"
$ws.Range("I294").Value = "1a nueva(o)"
$ws.Range("J294").Value = 400
$ws.Range("K294").Value = 8000
$ws.Range("L294").Value = 8500
$ws.Range("M294").Value = 8250
$ws.Range("N294").Value = "$/malla 18 kilos"
$ws.Range("O294").Value = "Región de O'Higgins"
$ws.Range("P294").Value = 458
$ws.Range("Q294").Value = 18
$ws.Range("R294").Value = "Hortaliza"
